# "Actualizar 02-04-2021 20-38-27" - automated availability-check update.
#
# 1) The 14-service block that was already stamped "44231.83819231426"
#    (rows 114-127, the run started at 2021-02-04 20:06:59) gets its
#    timestamp nudged to "44231.83819231481" (same instant, re-serialised).
# 2) A brand-new run stamped "44231.8595949117" (2021-02-04 20:37:49) is
#    appended as rows 128-141, cycling through the same 14 services/URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# --- 1) bump the timestamp already stored in rows 114-127 ---------------
for ($r = 114; $r -le 127; $r++) {
    $ws.Range("D$r").Value = 44231.83819231481
}

# --- 2) append the new availability-check run as rows 128-141 -----------
# Each tuple is (service name, URL) - identical cycle used throughout the
# sheet (Odoo, Blackbox, PowerBI, Dropbox, Odoo, GEE, UtilidadesOdoo,
# Filtros Dashboard, MapStore, GeoServer, Tomcat, Shiny, Github, EZ Exporter).
$services = @(
    @("Odoo", "https://www.dataintelligence-group.com/"),
    @("Blackbox", "https://serviciodashboard.azurewebsites.net/"),
    @("PowerBI", "https://powerbi.microsoft.com/es-es/"),
    @("Dropbox", "https://www.dropbox.com/"),
    @("Odoo", "https://dataintelligence.store/"),
    @("GEE", "https://app-data-i.users.earthengine.app/"),
    @("UtilidadesOdoo", "https://odooutil.azurewebsites.net/"),
    @("Filtros Dashboard", "https://filtradordashboard.azurewebsites.net/"),
    @("MapStore", "https://ide.dataintelligence-group.com/mapstore/#/"),
    @("GeoServer", "https://ide.dataintelligence-group.com/geoserver/web/?0"),
    @("Tomcat", "https://ide.dataintelligence-group.com/"),
    @("Shiny", "https://rpubs.com/dataintelligence/"),
    @("Github", "https://github.com/Sud-Austral/"),
    @("EZ Exporter", "https://ezexporter.highviewapps.com/exports/export-profile/")
)

$newTimestamp = 44231.8595949117
$startRow = 128

for ($i = 0; $i -lt $services.Count; $i++) {
    $row = $startRow + $i
    $name = $services[$i][0]
    $url = $services[$i][1]

    $ws.Range("A$row").Value = $name
    $ws.Range("C$row").Value = "Disponible"
    $ws.Range("D$row").Value = $newTimestamp
    $ws.Range("D$row").NumberFormat = $dateFormat

    $ws.Range("B$row").Value = $url

    $hashIndex = $url.IndexOf("#")
    if ($hashIndex -ge 0) {
        $address = $url.Substring(0, $hashIndex)
        $subAddress = $url.Substring($hashIndex + 1)
        $ws.Hyperlinks.Add($ws.Range("B$row"), $address, $subAddress)
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $url)
    }
    # Hyperlinks.Add stamps its own ad-hoc style on the cell - put the
    # sheet's normal "Hyperlink" style (shared by every other link cell)
    # back so the new cells stay consistent with B2:B127.
    $ws.Range("B$row").Style = "Hyperlink"
}

Write-Output "done"
